# Workbook/worksheet handles
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Full Name" column (A) for rows 2-13 -- it is being dropped in
# favor of relying solely on the separate First Name / Last Name columns,
# guarding downstream code against crashes when the combined-name cell is
# missing.
$ws.Range("A2").Value = $null
$ws.Range("A3").Value = $null
$ws.Range("A4").Value = $null
$ws.Range("A5").Value = $null
$ws.Range("A6").Value = $null
$ws.Range("A10").Value = $null
$ws.Range("A11").Value = $null
$ws.Range("A12").Value = $null
$ws.Range("A13").Value = $null

# Update First Name (D) / Last Name (E) values
$ws.Range("D2").Value = "Alexis"
$ws.Range("E2").Value = "Bernazzani"

$ws.Range("D3").Value = "Anna"
$ws.Range("E3").Value = "Thomas"

$ws.Range("D4").Value = "Ashley"

$ws.Range("D5").Value = "Claire"
$ws.Range("E5").Value = $null

$ws.Range("D6").Value = "Clara"

$ws.Range("D10").Value = "Jacky"

$ws.Range("D11").Value = "Juia"

$ws.Range("D12").Value = "Kayla"

$ws.Range("D13").Value = "Keala"

# Reposition the active cell selection to match the edited workbook.
[void]$ws.Range("F7").Select()
